# Saldo_guide.xlsx update: refresh the client-balance extract from
# 2024-10-25 to the 2024-10-28 pull (new report name + new "as of" dates),
# and correct a handful of balances that changed between runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to match the newer report extraction timestamp.
$ws.Name = "IClientBalance-20241028-094905-"

# Column G holds the "as of" date for every data row (2..274) as a date
# serial; the whole extract moved from 45590 (2024-10-25) to 45593
# (2024-10-28).
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45593
}

# A handful of balances (columns D/E/H) were recalculated between the two
# extracts; apply the corrected figures for those specific rows.
$ws.Range("E51").Value = 550.42999999999995
$ws.Range("H51").Value = 550.42999999999995

$ws.Range("E52").Value = 6329.16
$ws.Range("H52").Value = 6329.16

$ws.Range("E105").Value = 731.66
$ws.Range("H105").Value = 731.66

$ws.Range("E109").Value = 227.99
$ws.Range("H109").Value = 227.99

$ws.Range("E112").Value = 22.61
$ws.Range("H112").Value = 22.61

$ws.Range("E118").Value = 828.83
$ws.Range("H118").Value = 828.83

$ws.Range("E143").Value = 12709.69
$ws.Range("H143").Value = 12709.69

$ws.Range("E232").Value = 55910.3
$ws.Range("H232").Value = 55910.3

$ws.Range("D264").Value = 0
$ws.Range("E264").Value = 939.23
